# Big Data.docx edit: "6 bileşen" -> "5 bileşen", typed in pieces so the
# resulting runs mirror the real edit history (separate runs for "5",
# " bil", and "eşen vardır. ", with the _GoBack bookmark now sitting
# between " bil" and "eşen vardır. " instead of in the trailing empty
# paragraph).

$d = $word.ActiveDocument

# --- Step 1: change the digit "6" -> "5" ------------------------------
# Track revisions while doing this single-character edit so that, once
# accepted, the new "5" lands in its own run (matching real Word's
# behaviour of keeping separately-typed text in its own run) instead of
# being silently merged back into the surrounding text.
$d.TrackRevisions = $true

# Anchor the search on the unique phrase so we land on the right "6",
# then narrow the range down to just that digit before retyping it.
$digit = $d.Content
$digit.Find.Execute("6 bileşen vardır.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$digit = $d.Range($digit.Start, $digit.Start + 1)
$digit.Text = "5"

$d.AcceptAllRevisions()
$d.TrackRevisions = $false

# --- Step 2: split "bileşen" into "bil" | "eşen" and drop the _GoBack --
# bookmark right at that split point. Adding a bookmark in the middle of
# a run naturally forces Word to split the run in two around it, and
# because bookmark names must be unique, re-adding "_GoBack" here
# removes it from its old location (the trailing empty paragraph).
$split = $d.Content
$split.Find.Execute("5 bil", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $d.Range($split.End, $split.End)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null
